# Daily cryptos price/volume refresh (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) for every coin row, plus the three
# ranking swaps among rows 34-44 (coin, link, price, volume move together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.558.73'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '3.117.55'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.49'
$ws.Range('E5').Value = '  +8.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '628.70'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.05'
$ws.Range('E7').Value = '  +9.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.354'
$ws.Range('E8').Value = '  -9.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '3.113.23'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.719'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').Value = '  +4.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.48'
$ws.Range('E13').Value = '  +5.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.62'
$ws.Range('E14').Value = '  +4.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000242'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('D16').Value = '90.300.05'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').Value = '3.687.96'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '3.148.64'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.68'
$ws.Range('E19').Value = '  -2.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.34'
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('E21').Value = '  -4.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '448.72'
$ws.Range('E22').Value = '  +3.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.65'
$ws.Range('E23').Value = '  +9.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.98'
$ws.Range('E24').Value = '  +2.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.10'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.81'
$ws.Range('E26').Value = '  +4.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.48'
$ws.Range('E27').Value = '  +2.93%  '
$ws.Range('D28').Value = '3.275.94'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.50'
$ws.Range('E30').Value = '  +6.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.159'
$ws.Range('E31').Value = '  -4.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.37'
$ws.Range('E32').Value = '  +16.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.200'
$ws.Range('E33').Value = '  +34.01%  '
$ws.Range('B34').Value = 'dogwifhat'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.81'
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.150'
$ws.Range('E35').Value = '  +5.85%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '510.18'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.07'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.92'
$ws.Range('E38').Value = '  +3.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.32'
$ws.Range('E39').Value = '  +3.01%  '
$ws.Range('B40').Value = 'Binance-PegBSC-USD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.797'
$ws.Range('E40').Value = '  -20.24%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.425'
$ws.Range('E41').Value = '  +12.53%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.18'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0857'
$ws.Range('E43').Value = '  +2.90%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.38'
$ws.Range('E45').Value = '  +41.81%  '
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.697'
$ws.Range('E47').Value = '  +13.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '149.40'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.59'
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.54'
$ws.Range('E50').Value = '  +8.91%  '
$ws.Range('E51').Value = '  +4.20%  '
